$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "popularity" column
$ws.Range("S1").Value = "popularity"
# Match the formatting of the other header cells (e.g. R1)
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122) # xlPasteFormats

$popularity = @(86,88,89,84,84,88,85,85,83,84,85,94,82,74,84,79,84,89,58,80,80,82,82,85,84,87,80,81,87,89,83,81,85,80,6,85,77,83,86,85,79,77,88,86,78,80,78,81,84,84,82,90,77,78,89,0,74,27,85,80,84,83,85,82,82,67,74,73,71,84,82,81,77,79,78,81,84,76,79,78,90,76,85,80,89,75,84,74,63,74,84,82,84,78,84,78,77,78,82,84)

for ($i = 0; $i -lt $popularity.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $popularity[$i]
}
